$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51 (ALC)
$ws.Range("H51").Value = 5483.4
$ws.Range("J51").Value = 4083.8333
$ws.Range("L51").Value = 4083.8333
$ws.Range("N51").Value = -5051.8333

# Row 98 (ALC)
$ws.Range("H98").Value = 1767
$ws.Range("I98").Value = 1767
$ws.Range("K98").Value = 1767
$ws.Range("M98").Value = -269

# Row 100 (ALC)
$ws.Range("H100").Value = 3127.389
$ws.Range("I100").Value = 3525.5386
$ws.Range("K100").Value = 3525.5386
$ws.Range("M100").Value = -2984.5386

# Row 101 (ALC)
$ws.Range("H101").Value = 1684.3334
$ws.Range("J101").Value = 1798.8334
$ws.Range("L101").Value = 5396.5002
$ws.Range("N101").Value = -8640.5002

# Row 122 (ALC)
$ws.Range("H122").Value = 1767
$ws.Range("I122").Value = 1767
$ws.Range("K122").Value = 5301
$ws.Range("M122").Value = -2851

# Row 137 (ALC)
$ws.Range("H137").Value = 8848.174000000001
$ws.Range("I137").Value = 9331.5
$ws.Range("J137").Value = 5626
$ws.Range("K137").Value = 27994.5
$ws.Range("L137").Value = 16878
$ws.Range("M137").Value = -25444.5
$ws.Range("N137").Value = -21978

$ws = $wb.Worksheets.Item("ARM")
# Row 46 (ARM)
$ws.Range("H46").Value = 14976
$ws.Range("J46").Value = 14976
$ws.Range("L46").Value = 14976
$ws.Range("N46").Value = -15614

# Row 102 (ARM)
$ws.Range("H102").Value = 4030.6667
$ws.Range("I102").Value = 2183.1667
$ws.Range("K102").Value = 2183.1667
$ws.Range("M102").Value = -561.1667000000002

# Row 110 (ARM)
$ws.Range("H110").Value = 2010.6342
$ws.Range("I110").Value = 1963.5588
$ws.Range("K110").Value = 1963.5588
$ws.Range("M110").Value = 81.44119999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 2203.963
$ws.Range("I20").Value = 1657.9412
$ws.Range("J20").Value = 3132.2
$ws.Range("K20").Value = 1657.9412
$ws.Range("L20").Value = 3132.2
$ws.Range("M20").Value = -1410.9412
$ws.Range("N20").Value = -3626.2

# Row 22 (BSM)
$ws.Range("H22").Value = 247.5
$ws.Range("I22").Value = 247.5
$ws.Range("K22").Value = 247.5
$ws.Range("M22").Value = -74.5

# Row 94 (BSM)
$ws.Range("H94").Value = 572.4
$ws.Range("I94").Value = 572.4
$ws.Range("K94").Value = 572.4
$ws.Range("M94").Value = -121.4

# Row 107 (BSM)
$ws.Range("H107").Value = 3601.8333
$ws.Range("I107").Value = 2544.6
$ws.Range("J107").Value = 8888
$ws.Range("K107").Value = 2544.6
$ws.Range("L107").Value = 8888
$ws.Range("M107").Value = -624.5999999999999
$ws.Range("N107").Value = -12728

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 2578.0513
$ws.Range("I16").Value = 1479.8846
$ws.Range("K16").Value = 1479.8846
$ws.Range("M16").Value = -1192.8846

# Row 45 (CRP)
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 58 (CRP)
$ws.Range("H58").Value = 4294.6
$ws.Range("I58").Value = 4717.75
$ws.Range("J58").Value = 2602
$ws.Range("K58").Value = 4717.75
$ws.Range("L58").Value = 2602
$ws.Range("M58").Value = -4514.75
$ws.Range("N58").Value = -3008

# Row 105 (CRP)
$ws.Range("H105").Value = 1461.75
$ws.Range("I105").Value = 1437.7693
$ws.Range("J105").Value = 1565.6666
$ws.Range("K105").Value = 1437.7693
$ws.Range("L105").Value = 1565.6666
$ws.Range("M105").Value = 309.2307000000001
$ws.Range("N105").Value = -5059.6666

# Row 113 (CRP)
$ws.Range("H113").Value = 2578.0513
$ws.Range("I113").Value = 1479.8846
$ws.Range("K113").Value = 1479.8846
$ws.Range("M113").Value = 690.1153999999999

# Row 136 (CRP)
$ws.Range("H136").Value = 4294.6
$ws.Range("I136").Value = 4717.75
$ws.Range("J136").Value = 2602
$ws.Range("K136").Value = 14153.25
$ws.Range("L136").Value = 7806
$ws.Range("M136").Value = -11603.25
$ws.Range("N136").Value = -12906

$ws = $wb.Worksheets.Item("CUL")
# Row 37 (CUL)
$ws.Range("H37").Value = 117999.4
$ws.Range("J37").Value = 117999.4
$ws.Range("L37").Value = 353998.2
$ws.Range("N37").Value = -354222.2

# Row 68 (CUL)
$ws.Range("H68").Value = 4252.5
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 15000
$ws.Range("M68").Value = -14189

# Row 71 (CUL)
$ws.Range("H71").Value = 4252.5
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 45000
$ws.Range("M71").Value = -40944

$ws = $wb.Worksheets.Item("GSM")
# Row 107 (GSM)
$ws.Range("H107").Value = 6814.625
$ws.Range("I107").Value = 680.625
$ws.Range("J107").Value = 12948.625
$ws.Range("K107").Value = 680.625
$ws.Range("L107").Value = 12948.625
$ws.Range("M107").Value = 1239.375
$ws.Range("N107").Value = -16788.625

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2371.6667
$ws.Range("I7").Value = 2371.6667
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2371.6667
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -2259.6667

# Row 40 (LTW)
$ws.Range("H40").Value = 2397.2195
$ws.Range("I40").Value = 1949.069
$ws.Range("J40").Value = 3480.25
$ws.Range("K40").Value = 1949.069
$ws.Range("L40").Value = 3480.25
$ws.Range("M40").Value = -1813.069
$ws.Range("N40").Value = -3752.25

# Row 55 (LTW)
$ws.Range("H55").Value = 331.89474
$ws.Range("I55").Value = 344.6
$ws.Range("K55").Value = 344.6
$ws.Range("M55").Value = -171.6

# Row 68 (LTW)
$ws.Range("H68").Value = 25644064
$ws.Range("I68").Value = 30305894
$ws.Range("K68").Value = 30305894
$ws.Range("M68").Value = -30305145

# Row 71 (LTW)
$ws.Range("H71").Value = 25644064
$ws.Range("I71").Value = 30305894
$ws.Range("K71").Value = 151529470
$ws.Range("M71").Value = -151525726

# Row 93 (LTW)
$ws.Range("H93").Value = 8696808
$ws.Range("I93").Value = 16668603
$ws.Range("J93").Value = 303.72726
$ws.Range("K93").Value = 16668603
$ws.Range("L93").Value = 303.72726
$ws.Range("M93").Value = -16667355
$ws.Range("N93").Value = -2799.72726

# Row 126 (LTW)
$ws.Range("H126").Value = 2371.6667
$ws.Range("I126").Value = 2371.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7115.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -4645.000100000001

# Row 136 (LTW)
$ws.Range("H136").Value = 4793.931
$ws.Range("I136").Value = 3039.1155
$ws.Range("J136").Value = 20002.334
$ws.Range("K136").Value = 9117.3465
$ws.Range("L136").Value = 60007.00199999999
$ws.Range("M136").Value = -6567.3465
$ws.Range("N136").Value = -65107.00199999999

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws.Range("H96").Value = 4634.533
$ws.Range("I96").Value = 5016.2856
$ws.Range("J96").Value = 4300.5
$ws.Range("K96").Value = 5016.2856
$ws.Range("L96").Value = 4300.5
$ws.Range("M96").Value = -3643.2856
$ws.Range("N96").Value = -7046.5
